# Datorama_Mapping.xlsx edit:
#  - Add a new worksheet "Creative_Delivery_S3_Mapper" after "Creative_Delivery_Mapper",
#    populated like Creative_Delivery_Mapper but with sourceColumn == destinationTableColumnName
#    (creative delivery sql -> S3 mapping: source column names now match the destination names).
#  - Make the new sheet the active sheet/tab.
#  - Clear the old selection/active-cell state left on Creative_Delivery_Mapper and select A1:E28 there.

$wb = $excel.ActiveWorkbook

$deliveryMapper = $wb.Worksheets.Item("Creative_Delivery_Mapper")

# --- Create the new sheet right after Creative_Delivery_Mapper -----------------------------
$newSheet = $wb.Worksheets.Add($null, $deliveryMapper)
$newSheet.Name = "Creative_Delivery_S3_Mapper"

# --- Header row -------------------------------------------------------------------------------
$newSheet.Cells.Item(1, 1).Value = "sourceColumn"
$newSheet.Cells.Item(1, 2).Value = "destinationTableColumnName"
$newSheet.Cells.Item(1, 3).Value = "destinationTableUniqueColumn"
$newSheet.Cells.Item(1, 4).Value = "destinationTableDataType"
$newSheet.Cells.Item(1, 5).Value = "destinationTableValidationStyle"

# --- Data rows (sourceColumn mirrors destinationTableColumnName for the S3 feed) -------------
$rows = @(
  @("Date", "Date", $true, "DATE", "MATCH"),
  @("BU_ID", "BU_ID", $false, "INT", "MATCH"),
  @("Campaign_ID", "Campaign_ID", $false, "INT", "MATCH"),
  @("Campaing_Name", "Campaing_Name", $false, "VARCHAR", "MATCH"),
  @("Campaign_Flightdate_Start", "Campaign_Flightdate_Start", $false, "VARCHAR", "IGNORE"),
  @("Campaign_Flightdate_End", "Campaign_Flightdate_End", $false, "VARCHAR", "IGNORE"),
  @("Account_Manager_ID", "Account_Manager_ID", $false, "INT", "MATCH"),
  @("Campaign_Status", "Campaign_Status", $false, "VARCHAR", "MATCH"),
  @("Advertiser_Source_ID", "Advertiser_Source_ID", $false, "INT", "MATCH"),
  @("Advertiser_Source_Name", "Advertiser_Source_Name", $false, "VARCHAR", "MATCH"),
  @("Campaign_Target_ID", "Campaign_Target_ID", $true, "INT", "MATCH"),
  @("Campaign_Target_Name", "Campaign_Target_Name", $false, "VARCHAR", "MATCH"),
  @("Campaign_Target_Flightdate_Start", "Campaign_Target_Flightdate_Start", $false, "VARCHAR", "IGNORE"),
  @("Campaign_Target_Flightdate_End", "Campaign_Target_Flightdate_End", $false, "VARCHAR", "IGNORE"),
  @("Campaign_Target_Status", "Campaign_Target_Status", $false, "VARCHAR", "MATCH"),
  @("Creative_ID", "Creative_ID", $true, "VARCHAR", "MATCH"),
  @("Creative_Name", "Creative_Name", $false, "VARCHAR", "MATCH"),
  @("Creative_Message_ID", "Creative_Message_ID", $false, "INT", "MATCH"),
  @("Creative_Message_Name", "Creative_Message_Name", $false, "VARCHAR", "MATCH"),
  @("Adserver_Placement_ID", "Adserver_Placement_ID", $false, "VARCHAR", "MATCH"),
  @("Adserver_Placement_Name", "Adserver_Placement_Name", $false, "VARCHAR", "MATCH"),
  @("Integration_ID", "Integration_ID", $true, "INT", "MATCH"),
  @("Integration_Name", "Integration_Name", $false, "VARCHAR", "MATCH"),
  @("Currency_Code", "Currency_Code", $false, "VARCHAR", "MATCH"),
  @("Impressions", "Impressions", $false, "VARCHAR", "SUBSTRING"),
  @("Clicks", "Clicks", $false, "VARCHAR", "SUBSTRING"),
  @("Cost", "Cost", $false, "DOUBLE", "MATCH")
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# --- Column widths (approximate best-fit sizing for the new sheet) ---------------------------
$newSheet.Columns.Item(1).ColumnWidth = 25
$newSheet.Columns.Item(2).ColumnWidth = 32.140625
$newSheet.Columns.Item(3).ColumnWidth = 30
$newSheet.Columns.Item(4).ColumnWidth = 24.7109375
$newSheet.Columns.Item(5).ColumnWidth = 30.28515625

# --- Selection bookkeeping -----------------------------------------------------------------
# Clear the stale single-cell selection on Creative_Delivery_Mapper (was A26) and select the
# whole used range there instead, then leave the new sheet selected/active with A2:A28 chosen.
$deliveryMapper.Range("A1:E28").Select()
$newSheet.Range("A2:A28").Select()
